$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="G2"; Value=0.1294427768909185},
    @{Cell="J2"; Value=0.136021330371216},
    @{Cell="S2"; Value=-1.385286998382891},
    @{Cell="U2"; Value=1.385286998382891},
    @{Cell="C3"; Value=53.07356500808555},
    @{Cell="G3"; Value=0.1839686096833563},
    @{Cell="J3"; Value=0.1678799073680652},
    @{Cell="S3"; Value=-1.724321152095986},
    @{Cell="U3"; Value=1.724321152095986},
    @{Cell="C4"; Value=44.45195924760561},
    @{Cell="G4"; Value=0.1960820334845552},
    @{Cell="J4"; Value=0.1637440222852313},
    @{Cell="S4"; Value=-1.685476297326107},
    @{Cell="U4"; Value=1.685476297326107},
    @{Cell="C5"; Value=36.02457776097508},
    @{Cell="G5"; Value=0.1425237157371368},
    @{Cell="J5"; Value=0.1389360615571403},
    @{Cell="S5"; Value=-1.437368731193258},
    @{Cell="U5"; Value=1.437368731193258},
    @{Cell="C6"; Value=28.83773410500879},
    @{Cell="D6"; Value=-0.02066547782790784},
    @{Cell="F6"; Value=0.001033273891395392},
    @{Cell="G6"; Value=0.1854154611275955},
    @{Cell="H6"; Value=0.00202480351757841},
    @{Cell="J6"; Value=0.1711594340041738},
    @{Cell="R6"; Value=0.02066547782790784},
    @{Cell="S6"; Value=-1.74688134317385},
    @{Cell="U6"; Value=1.74688134317385},
    @{Cell="C7"; Value=20},
    @{Cell="D7"; Value=9.570888597260158},
    @{Cell="F7"; Value=0.4785444298630079},
    @{Cell="G7"; Value=-1.214500948661576},
    @{Cell="J7"; Value=0},
    @{Cell="P7"; Value=9.570888597260158},
    @{Cell="S7"; Value=-0},
    @{Cell="U7"; Value=0},
    @{Cell="C8"; Value=67.8544429863008},
    @{Cell="G8"; Value=0.09576465376964684},
    @{Cell="J8"; Value=0.1544808340706862},
    @{Cell="S8"; Value=-1.482399328957741},
    @{Cell="U8"; Value=1.482399328957741},
    @{Cell="C9"; Value=60.44244634151208},
    @{Cell="G9"; Value=0.1840019904277628},
    @{Cell="J9"; Value=0.1728313310115678},
    @{Cell="S9"; Value=-1.462070307178478},
    @{Cell="T9"; Value=-0},
    @{Cell="U9"; Value=1.462070307178478},
    @{Cell="W9"; Value=0},
    @{Cell="C10"; Value=53.13209480561969},
    @{Cell="G10"; Value=0.2625408050941181},
    @{Cell="J10"; Value=0.2350345021225493},
    @{Cell="S10"; Value=-1.620927600845167},
    @{Cell="T10"; Value=-0},
    @{Cell="U10"; Value=1.620927600845167},
    @{Cell="W10"; Value=0},
    @{Cell="C11"; Value=45.02745680139385},
    @{Cell="D11"; Value=0},
    @{Cell="F11"; Value=0},
    @{Cell="G11"; Value=0.2899672715031826},
    @{Cell="I11"; Value=0},
    @{Cell="J11"; Value=0.2513397664596594},
    @{Cell="O11"; Value=0},
    @{Cell="Q11"; Value=0},
    @{Cell="S11"; Value=-1.942047337812235},
    @{Cell="T11"; Value=-0},
    @{Cell="U11"; Value=1.942047337812235},
    @{Cell="W11"; Value=0},
    @{Cell="C12"; Value=35.31722011233268},
    @{Cell="D12"; Value=0},
    @{Cell="F12"; Value=0},
    @{Cell="G12"; Value=0.226580851315367},
    @{Cell="H12"; Value=0},
    @{Cell="I12"; Value=0},
    @{Cell="J12"; Value=0.1902380980480696},
    @{Cell="Q12"; Value=0},
    @{Cell="R12"; Value=0},
    @{Cell="S12"; Value=-1.886908332157008},
    @{Cell="T12"; Value=-0},
    @{Cell="U12"; Value=1.886908332157008},
    @{Cell="W12"; Value=0},
    @{Cell="C13"; Value=25.88267845154764},
    @{Cell="D13"; Value=-0.04914236872111433},
    @{Cell="F13"; Value=0.002457118436055716},
    @{Cell="G13"; Value=0.04379275037099467},
    @{Cell="H13"; Value=0.004625279744031281},
    @{Cell="I13"; Value=0},
    @{Cell="J13"; Value=0.1061102594279014},
    @{Cell="O13"; Value=0},
    @{Cell="Q13"; Value=0},
    @{Cell="R13"; Value=0.04914236872111433},
    @{Cell="S13"; Value=-1.127393321588413},
    @{Cell="T13"; Value=-0},
    @{Cell="U13"; Value=1.127393321588413},
    @{Cell="W13"; Value=0},
    @{Cell="C14"; Value=20},
    @{Cell="D14"; Value=3.217214041083077},
    @{Cell="F14"; Value=0.1608607020541538},
    @{Cell="G14"; Value=-0.4900704957521546},
    @{Cell="J14"; Value=0},
    @{Cell="O14"; Value=0},
    @{Cell="P14"; Value=3.217214041083077},
    @{Cell="T14"; Value=-0},
    @{Cell="U14"; Value=0},
    @{Cell="W14"; Value=0},
    @{Cell="C15"; Value=36.08607020541539},
    @{Cell="D15"; Value=0},
    @{Cell="F15"; Value=0},
    @{Cell="G15"; Value=0.1047120477578028},
    @{Cell="I15"; Value=0},
    @{Cell="J15"; Value=0.1258916198990845},
    @{Cell="O15"; Value=0},
    @{Cell="Q15"; Value=0},
    @{Cell="S15"; Value=-1.470696494148183},
    @{Cell="T15"; Value=-0},
    @{Cell="U15"; Value=1.470696494148183},
    @{Cell="W15"; Value=0},
    @{Cell="C16"; Value=28.73258773467447},
    @{Cell="D16"; Value=-0.08982950347167562},
    @{Cell="F16"; Value=0.004491475173583781},
    @{Cell="G16"; Value=0.1918553320931964},
    @{Cell="H16"; Value=0.008461939227031843},
    @{Cell="I16"; Value=0},
    @{Cell="J16"; Value=0.1560600136942353},
    @{Cell="O16"; Value=0},
    @{Cell="Q16"; Value=0},
    @{Cell="R16"; Value=0.08982950347167562},
    @{Cell="S16"; Value=-1.656688043463219},
    @{Cell="T16"; Value=-0},
    @{Cell="U16"; Value=1.656688043463219},
    @{Cell="W16"; Value=0},
    @{Cell="C17"; Value=20},
    @{Cell="D17"; Value=9.570888597260158},
    @{Cell="F17"; Value=0.4785444298630079},
    @{Cell="G17"; Value=-1.316801300087658},
    @{Cell="I17"; Value=0},
    @{Cell="J17"; Value=0},
    @{Cell="P17"; Value=9.570888597260158},
    @{Cell="Q17"; Value=0},
    @{Cell="T17"; Value=-0},
    @{Cell="U17"; Value=0},
    @{Cell="W17"; Value=0},
    @{Cell="C18"; Value=67.8544429863008},
    @{Cell="G18"; Value=-0.008694395993434634},
    @{Cell="I18"; Value=0},
    @{Cell="J18"; Value=0.09112933500273557},
    @{Cell="Q18"; Value=0},
    @{Cell="S18"; Value=-0.9232962006356187},
    @{Cell="T18"; Value=-0},
    @{Cell="U18"; Value=0.9232962006356187},
    @{Cell="W18"; Value=0},
    @{Cell="C19"; Value=63.23796198312269},
    @{Cell="G19"; Value=-0.1562102987179688},
    @{Cell="I19"; Value=0},
    @{Cell="J19"; Value=0.02285737553417968},
    @{Cell="Q19"; Value=0},
    @{Cell="S19"; Value=-0.2229987856993139},
    @{Cell="T19"; Value=-0},
    @{Cell="U19"; Value=0.2229987856993139},
    @{Cell="W19"; Value=0},
    @{Cell="C20"; Value=62.12296805462613},
    @{Cell="G20"; Value=0.2366472498453367},
    @{Cell="I20"; Value=0},
    @{Cell="J20"; Value=0.2061933707688903},
    @{Cell="Q20"; Value=0},
    @{Cell="S20"; Value=-1.911675975977103},
    @{Cell="T20"; Value=-0},
    @{Cell="U20"; Value=1.911675975977103},
    @{Cell="W20"; Value=0},
    @{Cell="C21"; Value=52.56458817474061},
    @{Cell="G21"; Value=0.280974756468153},
    @{Cell="J21"; Value=0.2368158151950638},
    @{Cell="S21"; Value=-1.568316656920952},
    @{Cell="T21"; Value=-0},
    @{Cell="U21"; Value=1.568316656920952},
    @{Cell="W21"; Value=0},
    @{Cell="C22"; Value=44.72300489013585},
    @{Cell="D22"; Value=0},
    @{Cell="F22"; Value=0},
    @{Cell="G22"; Value=0.3500536064002543},
    @{Cell="H22"; Value=0},
    @{Cell="J22"; Value=0.2987775943334393},
    @{Cell="R22"; Value=0},
    @{Cell="S22"; Value=-1.993844473362958},
    @{Cell="T22"; Value=-0},
    @{Cell="U22"; Value=1.993844473362958},
    @{Cell="W22"; Value=0},
    @{Cell="C23"; Value=34.75378252332106},
    @{Cell="G23"; Value=0.267690361013885},
    @{Cell="J23"; Value=0.2429292504224521},
    @{Cell="S23"; Value=-1.815614726625202},
    @{Cell="U23"; Value=1.815614726625202},
    @{Cell="C24"; Value=25.67570889019505},
    @{Cell="G24"; Value=-0.01069033521125611},
    @{Cell="J24"; Value=0.09191849366197663},
    @{Cell="S24"; Value=-0.885108268290579},
    @{Cell="U24"; Value=0.885108268290579},
    @{Cell="C25"; Value=21.25016754874216},
    @{Cell="G25"; Value=-0.1454136968490627},
    @{Cell="J25"; Value=0.02043495964622389},
    @{Cell="S25"; Value=-0.2128641629814988},
    @{Cell="U25"; Value=0.2128641629814988},
    @{Cell="C26"; Value=20.18584673383467},
    @{Cell="D26"; Value=0},
    @{Cell="F26"; Value=0},
    @{Cell="G26"; Value=-0.199987694159988},
    @{Cell="H26"; Value=0},
    @{Cell="J26"; Value=0.002120294100004973},
    @{Cell="R26"; Value=0},
    @{Cell="S26"; Value=-0.02059537736770251},
    @{Cell="U26"; Value=0.02059537736770251},
    @{Cell="C27"; Value=20.08286984699615}
)

foreach ($chg in $changes) {
    $ws.Range($chg.Cell).Value = $chg.Value
}